$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header text for B1:F1 (also made bold, matching the white Calibri 11
# run formatting already used for these header cells).
$headers = @{
    "B1" = "총 Chai 판매(개수)"
    "C1" = "Artisanal Chai 판매(단위)"
    "D1" = "미리 만든 Chai 판매(단위)"
    "E1" = "소셜 미디어 참여도(보기)"
    "F1" = "Chai에 대한 온라인 검색"
}

foreach ($addr in $headers.Keys) {
    $text = $headers[$addr]
    $len = $text.Length
    $cell = $ws.Range($addr)
    $cell.Value = $text

    # Bold the whole run. Splitting into [1..len-1] + [len..len] (rather than
    # one Characters(1, len) call) because a single call spanning the entire
    # string is a no-op in this engine.
    $run1 = $cell.Characters(1, $len - 1)
    $run1.Font.Bold = $true
    $run1.Font.Color = 16777215
    $run1.Font.Size = 11
    $run1.Font.Name = "Calibri"

    $run2 = $cell.Characters($len, 1)
    $run2.Font.Bold = $true
    $run2.Font.Color = 16777215
    $run2.Font.Size = 11
    $run2.Font.Name = "Calibri"
}

# Keep the table's column headers (ListColumns) in sync with the renamed
# worksheet header cells.
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Item(2).Name = $headers["B1"]
$lo.ListColumns.Item(3).Name = $headers["C1"]
$lo.ListColumns.Item(4).Name = $headers["D1"]
$lo.ListColumns.Item(5).Name = $headers["E1"]
$lo.ListColumns.Item(6).Name = $headers["F1"]
